$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 21.31
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value2 = "'78.9%"
$ws.Cells.Item(2, 7).Style = "Normal"

$ws.Cells.Item(3, 1).Value = 89.59
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value2 = "'95.3%"
$ws.Cells.Item(3, 7).Style = "Normal"

$ws.Cells.Item(4, 1).Value = 9.35
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value2 = "'34.6%"
$ws.Cells.Item(4, 7).Style = "Normal"

$ws.Cells.Item(5, 1).Value = 68.04
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value2 = "'81.0%"
$ws.Cells.Item(5, 7).Style = "Normal"

$ws.Cells.Item(6, 1).Value = 82.91
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value2 = "'98.7%"
$ws.Cells.Item(6, 7).Style = "Normal"

$ws.Cells.Item(7, 1).Value = 1767.45
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 4
$ws.Cells.Item(7, 5).Value = 23
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value2 = "'99.9%"
$ws.Cells.Item(7, 7).Style = "Normal"

$ws.Cells.Item(8, 1).Value = 5.38
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value2 = "'19.9%"
$ws.Cells.Item(8, 7).Style = "Normal"

$ws.Cells.Item(9, 1).Value = 660.56
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value2 = "'98.6%"
$ws.Cells.Item(9, 7).Style = "Normal"

$ws.Cells.Item(10, 1).Value = 20.06
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value2 = "'74.3%"
$ws.Cells.Item(10, 7).Style = "Normal"

$ws.Cells.Item(11, 1).Value = 0.27
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value2 = "'1.0%"
$ws.Cells.Item(11, 7).Style = "Normal"

$ws.Cells.Item(12, 1).Value = 3.42
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value2 = "'12.7%"
$ws.Cells.Item(12, 7).Style = "Normal"

$ws.Cells.Item(13, 1).Value = 0.17
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value2 = "'0.6%"
$ws.Cells.Item(13, 7).Style = "Normal"

$ws.Cells.Item(14, 1).Value = 7.75
$ws.Cells.Item(14, 2).Value = 0
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value2 = "'28.7%"
$ws.Cells.Item(14, 7).Style = "Normal"

$ws.Cells.Item(15, 1).Value = 4.07
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value2 = "'15.1%"
$ws.Cells.Item(15, 7).Style = "Normal"

$ws.Cells.Item(16, 1).Value = 42.32
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value2 = "'74.2%"
$ws.Cells.Item(16, 7).Style = "Normal"

$ws.Cells.Item(17, 1).Value = 97.05
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 2
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value2 = "'85.1%"
$ws.Cells.Item(17, 7).Style = "Normal"

$ws.Cells.Item(18, 1).Value = 70.17
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value2 = "'83.5%"
$ws.Cells.Item(18, 7).Style = "Normal"

$ws.Cells.Item(19, 1).Value = 2.38
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value2 = "'8.8%"
$ws.Cells.Item(19, 7).Style = "Normal"

$ws.Cells.Item(20, 1).Value = 5.51
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value2 = "'20.4%"
$ws.Cells.Item(20, 7).Style = "Normal"

$ws.Cells.Item(21, 1).Value = 0.29
$ws.Cells.Item(21, 2).Value = 0
$ws.Cells.Item(21, 3).Value = 1
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value2 = "'1.1%"
$ws.Cells.Item(21, 7).Style = "Normal"

$ws.Cells.Item(22, 1).Value = 272.09
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(22, 4).Value = 2
$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value2 = "'98.9%"
$ws.Cells.Item(22, 7).Style = "Normal"

$ws.Cells.Item(23, 1).Value = 11.37
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value2 = "'42.1%"
$ws.Cells.Item(23, 7).Style = "Normal"

$ws.Cells.Item(24, 1).Value = 66.61
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 1
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value2 = "'99.4%"
$ws.Cells.Item(24, 7).Style = "Normal"

$ws.Cells.Item(25, 1).Value = 17.03
$ws.Cells.Item(25, 2).Value = 0
$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value2 = "'63.1%"
$ws.Cells.Item(25, 7).Style = "Normal"

$ws.Cells.Item(26, 1).Value = 104.52
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 2
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value2 = "'91.7%"
$ws.Cells.Item(26, 7).Style = "Normal"

$ws.Cells.Item(27, 1).Value = 213.02
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = 2
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value2 = "'97.7%"
$ws.Cells.Item(27, 7).Style = "Normal"

$ws.Cells.Item(28, 1).Value = 27.75
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value2 = "'48.7%"
$ws.Cells.Item(28, 7).Style = "Normal"

$ws.Cells.Item(29, 1).Value = 3.46
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value2 = "'12.8%"
$ws.Cells.Item(29, 7).Style = "Normal"

$ws.Cells.Item(30, 1).Value = 217.04
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(30, 5).Value = 2
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value2 = "'99.6%"
$ws.Cells.Item(30, 7).Style = "Normal"

$ws.Cells.Item(31, 1).Value = 19.24
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value2 = "'71.3%"
$ws.Cells.Item(31, 7).Style = "Normal"

$ws.Range("A2:A31").Select()
